$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows that were fully removed from the dataset (RM 232 and SC 92),
# in descending order so row indices of earlier rows are unaffected.
$ws.Rows(28).Delete()
$ws.Rows(26).Delete()

# Apply the individual cell edits (value corrections / new-missing markers)
# that differ between the original and the target data set.
$ws.Range("E5").ClearContents()
$ws.Range("C6").Value = 15.1
$ws.Range("C8").ClearContents()
$ws.Range("E11").Value = -7.9
$ws.Range("C19").Value = 13.2
$ws.Range("E19").ClearContents()
$ws.Range("C21").ClearContents()
$ws.Range("C23").Value = 12.2
$ws.Range("E23").Value = -7.0
$ws.Range("E25").Value = -7.1
$ws.Range("B26").ClearContents()
$ws.Range("B27").Value = -20.4
$ws.Range("C27").ClearContents()
$ws.Range("E27").ClearContents()
$ws.Range("B29").ClearContents()
$ws.Range("C29").Value = 11.2
$ws.Range("E29").ClearContents()
$ws.Range("E30").Value = -5.7
$ws.Range("E33").Value = -10.7
